$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 30.07831966666667
$ws.Range("H2").Value = 90.234959
$ws.Range("I2").Value = 0.2269842729019557
$ws.Range("J2").Value = 0.2269842729019557
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 47.57896333333333
$ws.Range("N2").Value = 142.73689
$ws.Range("O2").Value = 0.450188452948237
$ws.Range("P2").Value = 0.4501884529482371
$ws.Range("Q2").Value = 1431.095268548612
$ws.Range("R2").Value = 12879.85741693751
$ws.Range("S2").Value = 0.1021856986613119
$ws.Range("T2").Value = 0.1021856986613119

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 30.07831966666667
$ws.Range("H3").Value = 90.234959
$ws.Range("I3").Value = 0.2269842729019557
$ws.Range("J3").Value = 0.2269842729019557
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.829723666666666
$ws.Range("N3").Value = 29.489171
$ws.Range("O3").Value = 0.09300808131111737
$ws.Range("P3").Value = 0.09300808131111739
$ws.Range("Q3").Value = 295.6615706809988
$ws.Range("R3").Value = 2660.954136128989
$ws.Range("S3").Value = 0.02111137171040995
$ws.Range("T3").Value = 0.02111137171040996

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 30.07831966666667
$ws.Range("H4").Value = 90.234959
$ws.Range("I4").Value = 0.2269842729019557
$ws.Range("J4").Value = 0.2269842729019557
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.933664
$ws.Range("N4").Value = 32.800992
$ws.Range("O4").Value = 0.1034534789405002
$ws.Range("P4").Value = 0.1034534789405003
$ws.Range("Q4").Value = 328.8662409199254
$ws.Range("R4").Value = 2959.796168279328
$ws.Range("S4").Value = 0.02348231269648724
$ws.Range("T4").Value = 0.02348231269648724

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 30.07831966666667
$ws.Range("H5").Value = 90.234959
$ws.Range("I5").Value = 0.2269842729019557
$ws.Range("J5").Value = 0.2269842729019557
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 37.34441866666666
$ws.Range("N5").Value = 112.033256
$ws.Range("O5").Value = 0.3533499868001453
$ws.Range("P5").Value = 0.3533499868001453
$ws.Range("Q5").Value = 1123.257362421834
$ws.Range("R5").Value = 10109.3162617965
$ws.Range("S5").Value = 0.08020488983374663
$ws.Range("T5").Value = 0.08020488983374664

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 31.96959266666667
$ws.Range("H6").Value = 95.90877800000001
$ws.Range("I6").Value = 0.2412566535243296
$ws.Range("J6").Value = 0.2412566535243296
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.57896333333333
$ws.Range("N6").Value = 142.73689
$ws.Range("O6").Value = 0.450188452948237
$ws.Range("P6").Value = 0.4501884529482371
$ws.Range("Q6").Value = 1521.080077268936
$ws.Range("R6").Value = 13689.72069542042
$ws.Range("S6").Value = 0.1086109596135868
$ws.Range("T6").Value = 0.1086109596135868

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 31.96959266666667
$ws.Range("H7").Value = 95.90877800000001
$ws.Range("I7").Value = 0.2412566535243296
$ws.Range("J7").Value = 0.2412566535243296
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.829723666666666
$ws.Range("N7").Value = 29.489171
$ws.Range("O7").Value = 0.09300808131111737
$ws.Range("P7").Value = 0.09300808131111739
$ws.Range("Q7").Value = 314.2522616492265
$ws.Range("R7").Value = 2828.270354843038
$ws.Range("S7").Value = 0.02243881844783892
$ws.Range("T7").Value = 0.02243881844783892

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 31.96959266666667
$ws.Range("H8").Value = 95.90877800000001
$ws.Range("I8").Value = 0.2412566535243296
$ws.Range("J8").Value = 0.2412566535243296
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.933664
$ws.Range("N8").Value = 32.800992
$ws.Range("O8").Value = 0.1034534789405002
$ws.Range("P8").Value = 0.1034534789405003
$ws.Range("Q8").Value = 349.5447844341974
$ws.Range("R8").Value = 3145.903059907776
$ws.Range("S8").Value = 0.0249588401246348
$ws.Range("T8").Value = 0.0249588401246348

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 31.96959266666667
$ws.Range("H9").Value = 95.90877800000001
$ws.Range("I9").Value = 0.2412566535243296
$ws.Range("J9").Value = 0.2412566535243296
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 37.34441866666666
$ws.Range("N9").Value = 112.033256
$ws.Range("O9").Value = 0.3533499868001453
$ws.Range("P9").Value = 0.3533499868001453
$ws.Range("Q9").Value = 1193.885853146797
$ws.Range("R9").Value = 10744.97267832117
$ws.Range("S9").Value = 0.0852480353382691
$ws.Range("T9").Value = 0.08524803533826911

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 55.83720666666667
$ws.Range("H10").Value = 167.51162
$ws.Range("I10").Value = 0.4213722008598541
$ws.Range("J10").Value = 0.4213722008598541
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 47.57896333333333
$ws.Range("N10").Value = 142.73689
$ws.Range("O10").Value = 0.450188452948237
$ws.Range("P10").Value = 0.4501884529482371
$ws.Range("Q10").Value = 2656.676408629089
$ws.Range("R10").Value = 23910.0876776618
$ws.Range("S10").Value = 0.1896968992204915
$ws.Range("T10").Value = 0.1896968992204915

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 55.83720666666667
$ws.Range("H11").Value = 167.51162
$ws.Range("I11").Value = 0.4213722008598541
$ws.Range("J11").Value = 0.4213722008598541
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.829723666666666
$ws.Range("N11").Value = 29.489171
$ws.Range("O11").Value = 0.09300808131111737
$ws.Range("P11").Value = 0.09300808131111739
$ws.Range("Q11").Value = 548.8643118518911
$ws.Range("R11").Value = 4939.77880666702
$ws.Range("S11").Value = 0.03919101991981779
$ws.Range("T11").Value = 0.0391910199198178

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 55.83720666666667
$ws.Range("H12").Value = 167.51162
$ws.Range("I12").Value = 0.4213722008598541
$ws.Range("J12").Value = 0.4213722008598541
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.933664
$ws.Range("N12").Value = 32.800992
$ws.Range("O12").Value = 0.1034534789405002
$ws.Range("P12").Value = 0.1034534789405003
$ws.Range("Q12").Value = 610.5052563918933
$ws.Range("R12").Value = 5494.54730752704
$ws.Range("S12").Value = 0.04359242010776716
$ws.Range("T12").Value = 0.04359242010776716

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 55.83720666666667
$ws.Range("H13").Value = 167.51162
$ws.Range("I13").Value = 0.4213722008598541
$ws.Range("J13").Value = 0.4213722008598541
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 37.34441866666666
$ws.Range("N13").Value = 112.033256
$ws.Range("O13").Value = 0.3533499868001453
$ws.Range("P13").Value = 0.3533499868001453
$ws.Range("Q13").Value = 2085.208022937191
$ws.Range("R13").Value = 18766.87220643472
$ws.Range("S13").Value = 0.1488918616117776
$ws.Range("T13").Value = 0.1488918616117776

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 14.62767266666667
$ws.Range("H14").Value = 43.883018
$ws.Range("I14").Value = 0.1103868727138606
$ws.Range("J14").Value = 0.1103868727138606
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 47.57896333333333
$ws.Range("N14").Value = 142.73689
$ws.Range("O14").Value = 0.450188452948237
$ws.Range("P14").Value = 0.4501884529482371
$ws.Range("Q14").Value = 695.9695014593354
$ws.Range("R14").Value = 6263.72551313402
$ws.Range("S14").Value = 0.04969489545284688
$ws.Range("T14").Value = 0.04969489545284689

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 14.62767266666667
$ws.Range("H15").Value = 43.883018
$ws.Range("I15").Value = 0.1103868727138606
$ws.Range("J15").Value = 0.1103868727138606
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.829723666666666
$ws.Range("N15").Value = 29.489171
$ws.Range("O15").Value = 0.09300808131111737
$ws.Range("P15").Value = 0.09300808131111739
$ws.Range("Q15").Value = 143.7859801997865
$ws.Range("R15").Value = 1294.073821798078
$ws.Range("S15").Value = 0.01026687123305071
$ws.Range("T15").Value = 0.01026687123305072

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 14.62767266666667
$ws.Range("H16").Value = 43.883018
$ws.Range("I16").Value = 0.1103868727138606
$ws.Range("J16").Value = 0.1103868727138606
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.933664
$ws.Range("N16").Value = 32.800992
$ws.Range("O16").Value = 0.1034534789405002
$ws.Range("P16").Value = 0.1034534789405003
$ws.Range("Q16").Value = 159.9340580393173
$ws.Range("R16").Value = 1439.406522353856
$ws.Range("S16").Value = 0.01141990601161106
$ws.Range("T16").Value = 0.01141990601161107

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 14.62767266666667
$ws.Range("H17").Value = 43.883018
$ws.Range("I17").Value = 0.1103868727138606
$ws.Range("J17").Value = 0.1103868727138606
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 37.34441866666666
$ws.Range("N17").Value = 112.033256
$ws.Range("O17").Value = 0.3533499868001453
$ws.Range("P17").Value = 0.3533499868001453
$ws.Range("Q17").Value = 546.2619321829565
$ws.Range("R17").Value = 4916.357389646608
$ws.Range("S17").Value = 0.03900520001635198
$ws.Range("T17").Value = 0.03900520001635199
